$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new icon names in column B for rows 32-35
$ws.Range("B32").Value = "download"
$ws.Range("B33").Value = "email"
$ws.Range("B34").Value = "details"
$ws.Range("B35").Value = "item-settings"

# Update the current selection to match the saved view state
$ws.Range("E32").Select()
